$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 11) of data, copying the formatting of an existing
# "normal" data row (row 9) so the new cells pick up the same styles.
$ws.Range("A9:C9").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = "Data10"
$ws.Range("B11").Value = "Noise  + numNeighbours"
$ws.Range("C11").Value = 20220408

# Update the selected cell, as recorded in the saved workbook.
$ws.Range("B12").Select()
